$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" column (R) is added to the right of the existing data table,
# mirroring the formatting of the preceding "2020" column (Q).

# Header cell: year 2021, formatted like Q4 (bold header row).
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# Data cell: 20.5 percent, formatted like Q5 (0.0 number format row).
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 20.5

# Clear clipboard marching ants / move the active selection to match the
# new state of the sheet (single cell S12, just past the new column).
$excel.CutCopyMode = $false
$ws.Range("S12").Select()
